# Generate Report for Handoff
# The 4 "Ready for handoff" rows (091d3e37, 33681839, 528ee6e2, 568544cf) in
# the zh-cn and de-de sheets get handed off: Priority flips from "low" to
# "ht", and the Latest Handoff Datetime timestamp is refreshed. The Overview
# sheet's "Latest HO Xliff Generate Date" column for those same source files
# picks up the new de-de generate timestamp as well (it shared the same
# underlying value as the de-de sheet's handoff time).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcnHandoffTime = "2016-08-21 02:37:12"
$dedeHandoffTime = "2016-08-21 02:37:17"

foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = $zhcnHandoffTime

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = $dedeHandoffTime

    $overview.Range("G$row").Value = $dedeHandoffTime
}
